$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure ion_id column (A) stays text, matching original formatting
$ws.Range("A2:A15").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2400"
$ws.Cells.Item(2, 2).Value = 594
$ws.Cells.Item(2, 3).Value = "[-12.920562520623207, 17.63073266670108, -7.267282247543335]"
$ws.Cells.Item(2, 4).Value = 23.03469256316377
$ws.Cells.Item(2, 5).Value = -8.033864768512299
$ws.Cells.Item(2, 6).Value = -0.3487723895807388
$ws.Cells.Item(2, 7).Value = 21.85826320213968
$ws.Cells.Item(2, 8).Value = -7.267282247543335
$ws.Cells.Item(2, 9).Value = "[3.486774444580078, -5.9720001220703125, -8.539749145507812]"

$ws.Cells.Item(3, 1).Value = "2399"
$ws.Cells.Item(3, 2).Value = 692
$ws.Cells.Item(3, 3).Value = "[8.804081737995148, 0.3854767978191376, -1.8715763092041016]"
$ws.Cells.Item(3, 4).Value = 9.009064629146058
$ws.Cells.Item(3, 5).Value = 2.104128715708714
$ws.Cells.Item(3, 6).Value = 0.2335568455021904
$ws.Cells.Item(3, 7).Value = 8.812516531102597
$ws.Cells.Item(3, 8).Value = -1.871576309204102
$ws.Cells.Item(3, 9).Value = "[0.3143463134765625, -1.2615242004394531, -8.960662841796875]"

$ws.Cells.Item(4, 1).Value = "2398"
$ws.Cells.Item(4, 2).Value = 565
$ws.Cells.Item(4, 3).Value = "[-0.7639022618532181, 4.1012390702962875, -19.930686593055725]"
$ws.Cells.Item(4, 4).Value = 20.36261222554711
$ws.Cells.Item(4, 5).Value = 19.13130058928689
$ws.Cells.Item(4, 6).Value = 0.9395307624276515
$ws.Cells.Item(4, 7).Value = 4.171775230928581
$ws.Cells.Item(4, 8).Value = -19.93068659305573
$ws.Cells.Item(4, 9).Value = "[3.2506446838378906, -1.1598472595214844, -17.24256134033203]"

$ws.Cells.Item(5, 1).Value = "2397"
$ws.Cells.Item(5, 2).Value = 969
$ws.Cells.Item(5, 3).Value = "[-0.6533946990966797, -10.039370059967041, -4.531984776258469]"
$ws.Cells.Item(5, 4).Value = 11.03425402308687
$ws.Cells.Item(5, 5).Value = 2.553530355577822
$ws.Cells.Item(5, 6).Value = 0.2314184855845346
$ws.Cells.Item(5, 7).Value = 10.06061011240224
$ws.Cells.Item(5, 8).Value = -4.531984776258469
$ws.Cells.Item(5, 9).Value = "[5.960777282714844, 2.729267120361328, -17.362136840820312]"

$ws.Cells.Item(6, 1).Value = "2276"
$ws.Cells.Item(6, 2).Value = 707
$ws.Cells.Item(6, 3).Value = "[0.7507225573062897, 4.773991197347641, -15.556918144226074]"
$ws.Cells.Item(6, 4).Value = 16.29025102497048
$ws.Cells.Item(6, 5).Value = 9.975358728624318
$ws.Cells.Item(6, 6).Value = 0.6123514434083066
$ws.Cells.Item(6, 7).Value = 4.832657272184864
$ws.Cells.Item(6, 8).Value = -15.55691814422607
$ws.Cells.Item(6, 9).Value = "[1.5978927612304688, -2.971912384033203, -4.38037109375]"

$ws.Cells.Item(7, 1).Value = "2520"
$ws.Cells.Item(7, 2).Value = 787
$ws.Cells.Item(7, 3).Value = "[-13.154578566551208, 3.3986346274614334, -9.698930323123932]"
$ws.Cells.Item(7, 4).Value = 16.69319933408097
$ws.Cells.Item(7, 5).Value = 3.148076982409002
$ws.Cells.Item(7, 6).Value = 0.1885844001144743
$ws.Cells.Item(7, 7).Value = 13.58652474308822
$ws.Cells.Item(7, 8).Value = -9.698930323123932
$ws.Cells.Item(7, 9).Value = "[7.958831787109375, 0.7162055969238281, -16.491683959960938]"

$ws.Cells.Item(8, 1).Value = "2223"
$ws.Cells.Item(8, 2).Value = 959
$ws.Cells.Item(8, 3).Value = "[4.85807591676712, 4.0646979212760925, -12.652371108531952]"
$ws.Cells.Item(8, 4).Value = 14.14938745926172
$ws.Cells.Item(8, 5).Value = 10.54082691116718
$ws.Cells.Item(8, 6).Value = 0.7449670129901985
$ws.Cells.Item(8, 7).Value = 6.334245874948247
$ws.Cells.Item(8, 8).Value = -12.65237110853195
$ws.Cells.Item(8, 9).Value = "[-1.7506599426269531, -5.281822204589844, -18.320533752441406]"

$ws.Cells.Item(9, 1).Value = "2391"
$ws.Cells.Item(9, 2).Value = 777
$ws.Cells.Item(9, 3).Value = "[2.6081194430589676, -4.2312188148498535, -9.195363759994507]"
$ws.Cells.Item(9, 4).Value = 10.45276108820162
$ws.Cells.Item(9, 5).Value = 8.878427138099886
$ws.Cells.Item(9, 6).Value = 0.8493858286038185
$ws.Cells.Item(9, 7).Value = 4.970462723771462
$ws.Cells.Item(9, 8).Value = -9.195363759994507
$ws.Cells.Item(9, 9).Value = "[0.9539413452148438, 1.3969001770019531, -13.710983276367188]"

$ws.Cells.Item(10, 1).Value = "2247"
$ws.Cells.Item(10, 2).Value = 794
$ws.Cells.Item(10, 3).Value = "[0.25486694276332855, 3.973426640033722, -6.730924963951111]"
$ws.Cells.Item(10, 4).Value = 7.820385367268257
$ws.Cells.Item(10, 5).Value = 3.182860451367143
$ws.Cells.Item(10, 6).Value = 0.4069953463788129
$ws.Cells.Item(10, 7).Value = 3.98159219687843
$ws.Cells.Item(10, 8).Value = -6.730924963951111
$ws.Cells.Item(10, 9).Value = "[2.0873985290527344, -2.63055419921875, -3.9116897583007812]"

$ws.Cells.Item(11, 1).Value = "2231"
$ws.Cells.Item(11, 2).Value = 806
$ws.Cells.Item(11, 3).Value = "[6.534557104110718, 3.7390658259391785, 4.73402202129364]"
$ws.Cells.Item(11, 4).Value = 8.893369119500393
$ws.Cells.Item(11, 5).Value = -4.397830002368352
$ws.Cells.Item(11, 6).Value = -0.4945066310949894
$ws.Cells.Item(11, 7).Value = 7.528681810090666
$ws.Cells.Item(11, 8).Value = 4.73402202129364
$ws.Cells.Item(11, 9).Value = "[0.74676513671875, 0.5153694152832031, -20.509132385253906]"

$ws.Cells.Item(12, 1).Value = "2381"
$ws.Cells.Item(12, 2).Value = 946
$ws.Cells.Item(12, 3).Value = "[9.601044237613678, -11.661227107048035, -1.7264618799090385]"
$ws.Cells.Item(12, 4).Value = 15.20345153961908
$ws.Cells.Item(12, 5).Value = -4.061560214754879
$ws.Cells.Item(12, 6).Value = -0.2671472464111686
$ws.Cells.Item(12, 7).Value = 15.10510735131554
$ws.Cells.Item(12, 8).Value = -1.726461879909039
$ws.Cells.Item(12, 9).Value = "[-0.46352386474609375, 6.322200775146484, -12.438140869140625]"

$ws.Cells.Item(13, 1).Value = "2343"
$ws.Cells.Item(13, 2).Value = 1186
$ws.Cells.Item(13, 3).Value = "[2.4428126215934753, 0.49578909622505307, -8.41463577747345]"
$ws.Cells.Item(13, 4).Value = 8.7760603689633
$ws.Cells.Item(13, 5).Value = 8.280676746972738
$ws.Cells.Item(13, 6).Value = 0.9435528470448431
$ws.Cells.Item(13, 7).Value = 2.492617165180414
$ws.Cells.Item(13, 8).Value = -8.41463577747345
$ws.Cells.Item(13, 9).Value = "[4.815853118896484, -3.690704345703125, -12.438545227050781]"

$ws.Cells.Item(14, 1).Value = "2230"
$ws.Cells.Item(14, 2).Value = 1225
$ws.Cells.Item(14, 3).Value = "[-0.9318812191486359, -2.220571478828788, -13.284280061721802]"
$ws.Cells.Item(14, 4).Value = 13.50079394174466
$ws.Cells.Item(14, 5).Value = 13.21380134611099
$ws.Cells.Item(14, 6).Value = 0.978742539374201
$ws.Cells.Item(14, 7).Value = 2.408181948937791
$ws.Cells.Item(14, 8).Value = -13.2842800617218
$ws.Cells.Item(14, 9).Value = "[1.012481689453125, -0.3558006286621094, -11.524864196777344]"

$ws.Cells.Item(15, 1).Value = "2515"
$ws.Cells.Item(15, 2).Value = 1219
$ws.Cells.Item(15, 3).Value = "[-7.849294036626816, 10.417173475027084, -12.752289831638336]"
$ws.Cells.Item(15, 4).Value = 18.24143130437794
$ws.Cells.Item(15, 5).Value = 5.720300019439035
$ws.Cells.Item(15, 6).Value = 0.3135883321867491
$ws.Cells.Item(15, 7).Value = 13.04334773293395
$ws.Cells.Item(15, 8).Value = -12.75228983163834
$ws.Cells.Item(15, 9).Value = "[9.918651580810547, -4.4640045166015625, -19.940650939941406]"
